# Applies the "New crime data collected" weekly CompStat update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header text updates (shared strings used in the title cells)
# ---------------------------------------------------------------------------
# "Volume 31   Number  51" -> "... 52"
$a8 = $ws.Range("A8").Value()
$ws.Range("A8").Value = $a8.Replace("51", "52")

# "Report Covering the Week  12/16/2024  Through  12/22/2024"
$c9 = $ws.Range("C9").Value()
$c9 = $c9.Replace("12/16/2024", "12/23/2024").Replace("12/22/2024", "12/29/2024")
$ws.Range("C9").Value = $c9

# ---------------------------------------------------------------------------
# Helper: write a text-typed cell (keeps/produces a shared-string cell with
# the "General" text style, s=13) by temporarily forcing a text number
# format, assigning the value, then restoring the original number format by
# copying it (format only) from a donor cell that already carries the
# desired style.
# ---------------------------------------------------------------------------
function Set-TextCell {
    param($ws, [string]$targetRef, [string]$donorRef, [string]$text)
    $ws.Range($targetRef).NumberFormat = "@"
    $ws.Range($targetRef).Value = $text
    $ws.Range($donorRef).Copy()
    $ws.Range($targetRef).PasteSpecial(-4122)
}

# Helper: write a numeric cell using a specific Excel number format string
# so the resulting style matches the workbook's existing numeric styles.
function Set-NumCell {
    param($ws, [string]$ref, [string]$numFmt, $value)
    $ws.Range($ref).NumberFormat = $numFmt
    $ws.Range($ref).Value = $value
}

$FMT_COUNT = "#,##0"
$FMT_PCT   = "#,##0.0;""-""#,##0.0"

# ---------------------------------------------------------------------------
# 2. Row 14
# ---------------------------------------------------------------------------
Set-NumCell $ws "N14" $FMT_PCT -84.615384615384

# ---------------------------------------------------------------------------
# 3. Row 15 (Rape)
# ---------------------------------------------------------------------------
Set-NumCell $ws "D15" $FMT_COUNT 1
Set-NumCell $ws "E15" $FMT_PCT -100
Set-NumCell $ws "F15" $FMT_COUNT 1
Set-NumCell $ws "G15" $FMT_COUNT 2
Set-NumCell $ws "H15" $FMT_PCT -50
Set-NumCell $ws "J15" $FMT_COUNT 19
Set-NumCell $ws "K15" $FMT_PCT -31.578947368421
Set-NumCell $ws "L15" $FMT_PCT -40.909090909090

# ---------------------------------------------------------------------------
# 4. Row 16 (Robbery)
# ---------------------------------------------------------------------------
Set-TextCell $ws "C16" "C15" "0"
Set-NumCell $ws "D16" $FMT_COUNT 2
Set-NumCell $ws "E16" $FMT_PCT -100
Set-NumCell $ws "F16" $FMT_COUNT 8
Set-NumCell $ws "G16" $FMT_COUNT 11
Set-NumCell $ws "H16" $FMT_PCT -27.272727272727
Set-NumCell $ws "I16" $FMT_COUNT 129
Set-NumCell $ws "J16" $FMT_COUNT 152
Set-NumCell $ws "K16" $FMT_PCT -15.131578947368
Set-NumCell $ws "L16" $FMT_PCT -15.686274509803
Set-NumCell $ws "M16" $FMT_PCT -63.142857142857
Set-NumCell $ws "N16" $FMT_PCT -89.486552567237

# ---------------------------------------------------------------------------
# 5. Row 17 (Fel. Assault)
# ---------------------------------------------------------------------------
Set-NumCell $ws "C17" $FMT_COUNT 2
Set-NumCell $ws "D17" $FMT_COUNT 5
Set-NumCell $ws "E17" $FMT_PCT -60
Set-NumCell $ws "F17" $FMT_COUNT 7
Set-NumCell $ws "G17" $FMT_COUNT 26
Set-NumCell $ws "H17" $FMT_PCT -73.076923076923
Set-NumCell $ws "I17" $FMT_COUNT 282
Set-NumCell $ws "J17" $FMT_COUNT 304
Set-NumCell $ws "K17" $FMT_PCT -7.236842105263
Set-NumCell $ws "L17" $FMT_PCT -15.820895522388
Set-NumCell $ws "M17" $FMT_PCT -13.496932515337
Set-NumCell $ws "N17" $FMT_PCT -66.784452296819

# ---------------------------------------------------------------------------
# 6. Row 18 (Burglary)
# ---------------------------------------------------------------------------
Set-TextCell $ws "C18" "C15" "0"
Set-NumCell $ws "D18" $FMT_COUNT 4
Set-NumCell $ws "E18" $FMT_PCT -100
Set-NumCell $ws "F18" $FMT_COUNT 4
Set-NumCell $ws "G18" $FMT_COUNT 11
Set-NumCell $ws "H18" $FMT_PCT -63.636363636363
Set-NumCell $ws "J18" $FMT_COUNT 139
Set-NumCell $ws "K18" $FMT_PCT -43.165467625899
Set-NumCell $ws "L18" $FMT_PCT -61.835748792270
Set-NumCell $ws "M18" $FMT_PCT -68.016194331983
Set-NumCell $ws "N18" $FMT_PCT -90.781796966161

# ---------------------------------------------------------------------------
# 7. Row 19 (Gr. Larceny)
# ---------------------------------------------------------------------------
Set-NumCell $ws "D19" $FMT_COUNT 3
Set-NumCell $ws "E19" $FMT_PCT 0
Set-NumCell $ws "F19" $FMT_COUNT 18
Set-NumCell $ws "G19" $FMT_COUNT 17
Set-NumCell $ws "H19" $FMT_PCT 5.882352941176
Set-NumCell $ws "I19" $FMT_COUNT 270
Set-NumCell $ws "J19" $FMT_COUNT 339
Set-NumCell $ws "K19" $FMT_PCT -20.353982300885
Set-NumCell $ws "L19" $FMT_PCT -37.354988399071
Set-NumCell $ws "M19" $FMT_PCT -25.824175824175
Set-NumCell $ws "N19" $FMT_PCT -29.503916449086

# ---------------------------------------------------------------------------
# 8. Row 20 (G.L.A.)
# ---------------------------------------------------------------------------
Set-NumCell $ws "C20" $FMT_COUNT 1
Set-TextCell $ws "D20" "C15" "0"
Set-TextCell $ws "E20" "E22" "***.*"
Set-NumCell $ws "F20" $FMT_COUNT 2
Set-NumCell $ws "H20" $FMT_PCT -75
Set-NumCell $ws "I20" $FMT_COUNT 88
Set-NumCell $ws "K20" $FMT_PCT -16.981132075471
Set-NumCell $ws "L20" $FMT_PCT -35.766423357664
Set-NumCell $ws "M20" $FMT_PCT 6.024096385542
Set-NumCell $ws "N20" $FMT_PCT -85.008517887563

# ---------------------------------------------------------------------------
# 9. Row 21 (TOTAL - bold)
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 6
$ws.Range("D21").Value = 15
$ws.Range("E21").Value = -60
$ws.Range("F21").Value = 40
$ws.Range("G21").Value = 76
$ws.Range("H21").Value = -47.368421052631
$ws.Range("I21").Value = 865
$ws.Range("J21").Value = 1065
$ws.Range("K21").Value = -18.779342723004
$ws.Range("L21").Value = -32.945736434108
$ws.Range("M21").Value = -38.695960311835
$ws.Range("N21").Value = -78.418163672654

# ---------------------------------------------------------------------------
# 10. Row 22 (Transit)
# ---------------------------------------------------------------------------
Set-NumCell $ws "F22" $FMT_COUNT 2
Set-TextCell $ws "G22" "C22" "0"
Set-TextCell $ws "H22" "E22" "***.*"

# ---------------------------------------------------------------------------
# 11. Row 23 (Housing)
# ---------------------------------------------------------------------------
Set-TextCell $ws "C23" "C15" "0"
Set-NumCell $ws "E23" $FMT_PCT -100
Set-NumCell $ws "J23" $FMT_COUNT 86
Set-NumCell $ws "K23" $FMT_PCT 0
Set-NumCell $ws "M23" $FMT_PCT 1.176470588235

# ---------------------------------------------------------------------------
# 12. Row 24 (Petit Larceny)
# ---------------------------------------------------------------------------
Set-NumCell $ws "C24" $FMT_COUNT 6
Set-NumCell $ws "D24" $FMT_COUNT 9
Set-NumCell $ws "E24" $FMT_PCT -33.333333333333
Set-NumCell $ws "F24" $FMT_COUNT 56
Set-NumCell $ws "G24" $FMT_COUNT 55
Set-NumCell $ws "H24" $FMT_PCT 1.818181818181
Set-NumCell $ws "I24" $FMT_COUNT 724
Set-NumCell $ws "J24" $FMT_COUNT 836
Set-NumCell $ws "K24" $FMT_PCT -13.397129186602
Set-NumCell $ws "L24" $FMT_PCT -17.820658342792
Set-NumCell $ws "M24" $FMT_PCT -9.612983770287

# ---------------------------------------------------------------------------
# 13. Row 25 (Retail Theft)
# ---------------------------------------------------------------------------
Set-TextCell $ws "C25" "C15" "0"
Set-NumCell $ws "D25" $FMT_COUNT 3
Set-NumCell $ws "E25" $FMT_PCT -100
Set-NumCell $ws "F25" $FMT_COUNT 8
Set-NumCell $ws "H25" $FMT_PCT 60
Set-NumCell $ws "I25" $FMT_COUNT 133
Set-NumCell $ws "J25" $FMT_COUNT 188
Set-NumCell $ws "K25" $FMT_PCT -29.255319148936
Set-NumCell $ws "L25" $FMT_PCT -48.846153846153

# ---------------------------------------------------------------------------
# 14. Row 26 (Misd. Assault)
# ---------------------------------------------------------------------------
Set-NumCell $ws "C26" $FMT_COUNT 8
Set-NumCell $ws "D26" $FMT_COUNT 9
Set-NumCell $ws "E26" $FMT_PCT -11.111111111111
Set-NumCell $ws "F26" $FMT_COUNT 32
Set-NumCell $ws "G26" $FMT_COUNT 32
Set-NumCell $ws "I26" $FMT_COUNT 417
Set-NumCell $ws "J26" $FMT_COUNT 494
Set-NumCell $ws "K26" $FMT_PCT -15.587044534413
Set-NumCell $ws "L26" $FMT_PCT 0.481927710843
Set-NumCell $ws "M26" $FMT_PCT -47.081218274111

# ---------------------------------------------------------------------------
# 15. Row 27 (UCR Rape*)
# ---------------------------------------------------------------------------
Set-NumCell $ws "D27" $FMT_COUNT 1
Set-NumCell $ws "E27" $FMT_PCT -100
Set-NumCell $ws "F27" $FMT_COUNT 1
Set-NumCell $ws "G27" $FMT_COUNT 2
Set-NumCell $ws "H27" $FMT_PCT -50
Set-NumCell $ws "J27" $FMT_COUNT 27
Set-NumCell $ws "K27" $FMT_PCT -29.629629629629
Set-NumCell $ws "L27" $FMT_PCT -36.666666666666

# ---------------------------------------------------------------------------
# 16. Row 28 (Other Sex Crimes)
# ---------------------------------------------------------------------------
Set-TextCell $ws "C28" "C15" "0"
Set-NumCell $ws "D28" $FMT_COUNT 1
Set-NumCell $ws "E28" $FMT_PCT -100
Set-NumCell $ws "F28" $FMT_COUNT 4
Set-NumCell $ws "H28" $FMT_PCT 300
Set-NumCell $ws "J28" $FMT_COUNT 30
Set-NumCell $ws "K28" $FMT_PCT 33.333333333333
Set-NumCell $ws "L28" $FMT_PCT 11.111111111111

# ---------------------------------------------------------------------------
# 17. Row 29 (Shooting Vic.)
# ---------------------------------------------------------------------------
Set-NumCell $ws "G29" $FMT_COUNT 3
Set-NumCell $ws "M29" $FMT_PCT -68.333333333333

# ---------------------------------------------------------------------------
# 18. Row 30 (Shooting Inc.)
# ---------------------------------------------------------------------------
Set-NumCell $ws "G30" $FMT_COUNT 2
Set-NumCell $ws "M30" $FMT_PCT -68.627450980392
